$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet already has four Longitude/Latitude column pairs (B:C, D:E, F:G, H:I).
# Add a fifth pair in columns J:K that duplicates the first pair (B:C) for every
# data row (header in row 1, data in rows 2-50).
for ($r = 1; $r -le 50; $r++) {
    $lon = $ws.Cells.Item($r, 2).Formula   # column B
    $lat = $ws.Cells.Item($r, 3).Formula   # column C
    $ws.Cells.Item($r, 10).Formula = $lon  # column J
    $ws.Cells.Item($r, 11).Formula = $lat  # column K
}

# Give the two new columns explicit widths (mirroring columns B and C).
$ws.Columns.Item(10).ColumnWidth = 19.666666666666668
$ws.Columns.Item(11).ColumnWidth = 16

# Update the view: scroll the window back to the top-left and select C7
# (previously the view was scrolled to row 32 with A50 selected).
$ws.Range("C7").Select() | Out-Null
